$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 44898
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("B9").Value = "vandaag ben ik de lessen gevolgd en nog wat onderzoek gedaan"

# Row 11 (date cell present but empty, keeps the date style)
$ws.Range("A11").NumberFormat = "d-mmm"

# Row 12
$ws.Range("A12").Value = 44901
$ws.Range("A12").NumberFormat = "d-mmm"
$ws.Range("B12").Value = "vandaag heb ik de lessen gevolgd."

# Row 13
$ws.Range("A13").Value = 44902
$ws.Range("A13").NumberFormat = "d-mmm"
$ws.Range("B13").Value = "vandaag heb ik de lessen gevolgd van ed en ben verder gegaan met onderzoek naar het project"

# Row 14
$ws.Range("A14").Value = 44903
$ws.Range("A14").NumberFormat = "d-mmm"
$ws.Range("B14").Value = "vandaag heb ik de informatie pagina van  onze kunstenaar gemaakt"

# Row 15
$ws.Range("A15").Value = 44904
$ws.Range("A15").NumberFormat = "d-mmm"
$ws.Range("B15").Value = "vandaag heb ik een begin gemaakt aan de form pagina "

# Row 16
$ws.Range("A16").Value = 44905
$ws.Range("A16").NumberFormat = "d-mmm"
$ws.Range("B16").Value = "vandaag heb ik de lessen gevolgd."

# Row 18
$ws.Range("A18").Value = 44908
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("B18").Value = "studiedag"

# Row 19
$ws.Range("A19").Value = 44909
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("B19").Value = "vrij"

# Row 20
$ws.Range("A20").Value = 44910
$ws.Range("A20").NumberFormat = "d-mmm"
$ws.Range("B20").Value = "vandaag heb ik met mijn team bedacht wat we nog meer nodig hadden qua spullen"

# Row 21
$ws.Range("A21").Value = 44911
$ws.Range("A21").NumberFormat = "d-mmm"
$ws.Range("B21").Value = "vandaag heb ik de spullen gehaald met elisa gehaald zoals hout en kippengaas enzo"

# Row 22
$ws.Range("A22").Value = 44912
$ws.Range("A22").NumberFormat = "d-mmm"
$ws.Range("B22").Value = "vandaag heb ik de lessen gevolgd."

# Row 24
$ws.Range("A24").Value = 44906
$ws.Range("A24").NumberFormat = "d-mmm"
$ws.Range("B24").Value = "vandaag heb ik onderzoek gedaan naar database en form"

# Row 25
$ws.Range("A25").Value = 44907
$ws.Range("A25").NumberFormat = "d-mmm"

# Sheet view: scroll position and active-cell selection
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B24").Select() | Out-Null
